$wb = $excel.ActiveWorkbook

# Fix typo in sheet name: "usssenate" -> "ussenate"
$wsSenate = $wb.Worksheets.Item("usssenate")
$wsSenate.Name = "ussenate"

# Update the remembered selection on the ussenate sheet
$wsSenate.Range("C39").Select() | Out-Null

# Update the remembered selection on the statehou68 sheet
$wsStatehou68 = $wb.Worksheets.Item("statehou68")
$wsStatehou68.Range("A9").Select() | Out-Null

# Make statehou68 the active (selected) tab
$wsStatehou68.Activate() | Out-Null
